$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: duplicate of row 2 (VieonDpoint / Dpoint@2021 with mailto hyperlink)
$ws.Range("A3").Value = "VieonDpoint"
$ws.Range("B3").Value = "Dpoint@2021"
[void]$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Dpoint@2021")
$ws.Range("B3").Style = "Hyperlink"

# Row 4: same pattern again
$ws.Range("A4").Value = "VieonDpoint"
$ws.Range("B4").Value = "Dpoint@2021"
[void]$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:Dpoint@2021")
$ws.Range("B4").Style = "Hyperlink"

# Move the selection to match the author's final cursor position
[void]$ws.Range("G4").Select()
